$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 1 is the sheet's newest-date-first header row. The 9 new dates are
#    all newer than anything currently on the sheet, so they belong at the
#    very front (left) of row 1. Manually shift the existing B1:V1 headers
#    9 columns to the right (K1:AE1), working right-to-left so values are
#    not clobbered before they are copied. This only touches row 1 - the
#    data rows below keep their existing column positions untouched.
for ($c = 22; $c -ge 2; $c--) {
    $srcValue = $ws.Cells.Item(1, $c).Value()
    $ws.Cells.Item(1, $c + 9).Value = $srcValue
}

# Now fill the freed-up B1:J1 cells with the new dates, newest first.
$newDates = @("Sep_08", "Aug_25", "Aug_04", "Jul_23", "Jul_17", "Jul_07", "Jun_30", "Jun_24", "Jun_16")
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $newDates[$i]
}

# 2) Data rows (2 and down): each row keeps its own existing cells untouched;
#    9 more "UN" cells are appended immediately after whatever that row's
#    current last populated column is (rows are "ragged" - some analysts
#    have fewer weeks of history than others, so the append point differs
#    per row).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $lastCol = $ws.Cells.Item($r, 256).End(-4159).Column
    if ($lastCol -lt 2) {
        continue
    }
    for ($i = 1; $i -le 9; $i++) {
        $ws.Cells.Item($r, $lastCol + $i).Value = "UN"
    }
}
